$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A3: store as plain text "13/2/2020" (not converted to a date serial)
$ws.Cells.Item(3, 1).Value2 = "13/2/2020"

# B3: the note text
$ws.Cells.Item(3, 2).Value = "Sửa lại các khóa ngoại, bỏ note hình ảnh hợp đồng, sửa tên folder db"

# C3: author name, same as C2
$ws.Cells.Item(3, 3).Value = "Nghĩa"
